# #4801_PreRegister - add a "Run Date" input row (row 3) above the
# existing header/detail rows, and move the active selection to E16.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 3 is a bit shorter than the default (matches the source row).
$ws.Rows.Item(3).RowHeight = 18

# A3 ("Run Date" label): reuse the bold/bordered header look (font +
# border) already used on row 5 by copying its formats, then override the
# fill to the existing yellow highlight and left/ (default) vertical align.
$ws.Range("A5").Copy() | Out-Null
$ws.Range("A3").PasteSpecial(-4122) | Out-Null
$ws.Range("A3").Value = "Run Date"
$ws.Range("A3").Interior.Color = 65535
$ws.Range("A3").HorizontalAlignment = -4131
$ws.Range("A3").VerticalAlignment = -4107

# B3 is the (empty) input cell next to the label - pale-yellow highlight,
# default font/border.
$ws.Range("B3").Interior.Color = 10092543

# Move / record the current selection, matching the saved view state.
$ws.Range("E16").Select() | Out-Null
